$p = $ppt.ActivePresentation

function Replace-ExactText($range, $old, $new) {
    $full = $range.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw ("Text not found: " + $old)
    }
    $start = $idx + 1
    $len = $old.Length
    $c = $range.Characters($start, $len)
    $c.Text = $new
}

# --- Slide 13: "Bert ... question-answering" paragraph -- collapse the
#     quoted expand_split_sentences function-name callout (3 runs) into a
#     single plain-prose run ---
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
Replace-ExactText $tr13 'Bert is a highly effective model for question-answering tasks, but it has a limitation of handling only 512 tokens at a time. This becomes problematic when dealing with longer documents. To address this limitation, I developed the "expand_split_sentences" function. This function splits and expands sentences within the document, creating paragraphs with fewer than 512 tokens. Each paragraph is then converted into a data frame. Since multiple data frames may contain the correct answer, we determine the best answer by identifying the data frame with the highest start score' 'Bert is a highly effective model for question-answering tasks, but it has a limitation of handling only 512 tokens at a time. This becomes problematic when dealing with longer documents. To address this limitation, I developed a method. The method splits and expands sentences within the document, creating paragraphs with fewer than 512 tokens. Each paragraph is then converted into a data frame. Since multiple data frames may contain the correct answer, we determine the best answer by identifying the data frame with the highest start score'

# --- Slide 6: "Life insurance is a legally enforceable contract ..." --
#     wording / grammar corrections scattered across the long paragraph ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange

Replace-ExactText $tr6 'Life insurance is a legally enforceable contract between two parties both of whom are legally qualified to contract. It is therefore, necessary that the terms and conditions of the agreement must be suitably documented in a manner that would make it clear that both parties to the contract are Ad- idem i.e., of the same mind. Ad-Idem means that both the parties understand the same thing in the same sense or agree on the same subject. There must be a consensus or Ad-Idem between the parties to the ' 'Life insurance is a legally enforceable contract between two parties both of whom are legally qualified to contract. It is therefore, necessary that the terms and conditions of the agreement must be suitably documented in a manner that would make it clear that both parties to the contract are Ad- idem i.e., of the same mind. Ad-Idem means that both the parties understand the same thing in the same sense or are of the same mind on the same subject. There must be consensus or Ad-Idem between the parties to the '
Replace-ExactText $tr6 ' is possible provided all the terms and conditions, rights and duties - privileges and obligations are properly documented in terms that can be interpreted in a court of law. Between two human beings, sometimes silence means acceptance. But as the insurer is a legal personality entitled to a contract verbal discussion between parties to the contract is not possible and hence there is a need for ' ' is possible provided all the terms and conditions, rights and duties - privileges and obligations are properly documented in terms which can be clearly interpreted in a court of law. Between two human beings sometime silence means an acceptance. But as the insurer is a legal personality entitled to contract verbal discussion between parties to the contract is not possible and hence there is a need for '
Replace-ExactText $tr6 ' is also a contract of utmost good faith and enforced only in the distant future. It is therefore necessary that the declarations made by both parties should be put in black and white for future reference. Any suppression, willful, and material shall make the contract void. The insured, therefore, must declare all that he knows about himself, his health, and his financial status in answering questions contained in the proposal form and other ancillary documents which may be required by the ' ' is also a contract of utmost good faith and enforced only in the distant future. It is therefore necessary that the declarations made by both the parties should be put in black and white for future reference. Any suppression, willful and material shall make the contract void. The insured, therefore, has a duty to declare all that he knows about himself, his health, his financial status in answering questions contained in the proposal form and other ancillary documents which may be required by the '
Replace-ExactText $tr6 ' is an important factor in deciding the quantum of premium against a policy. The document proving the age, i.e. age proof must be reliable and the insured has to undertake its ' ' is an important factor in deciding the quantum of premium against a policy. The document proving the age, i.e. age proof must be reliable and the insured has to undertake as to its '
Replace-ExactText $tr6 '-standard age proofs are those which are comparatively less reliable and therefore the insurer accepts them with a pinch of salt. In other words, the insurer takes certain precautions before accepting such age proofs as ' '-standard age proofs are those which are comparatively less reliable and therefore the insurer accepts them with a pinch of salt. In other words the insurer takes certain precautions before accepting such age proofs as '
Replace-ExactText $tr6 ' of income is the document that may become necessary whenever the sum proposed is very high. Normally a sum proposed which is seven to eight times of the declared income is acceptable for insurance. But proposals do come to the insurer when the known source of income of the proposer is much less compared to the amount of insurance desired. A service holder normally does not face this problem as his sources of income are ' ' of income is the document may become necessary whenever the sum proposed is very high. Normally a sum proposed which is seven to eight times of the declared income is acceptable for insurance. But proposals do come to the insurer when the known source of income of the proposer is much less compared to the amount of insurance desired. A service holder normally does not face this problem as his sources of income are '
Replace-ExactText $tr6 ' needs to be explained that the policy is a valuable document and needs to be kept in safe custody and in the knowledge of the close relatives.' ' needs to be explained that the policy is a valuable document and needs to be kept in safe custody and in the knowledge of the close relatives..'
